# This script applies the "new naming convention" edits described in the
# commit diff across the three worksheets of the workbook:
#   - SV_calls
#   - CNV_and_Aneuploidy_calls
#   - CNV_metrics

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: SV_calls
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("SV_calls")

# Header renames
$ws1.Range("B1").Value = "Treated Sample Name"
$ws1.Range("E1").Value = "Start Chromosome"
$ws1.Range("F1").Value = "End Chromosome"
$ws1.Range("G1").Value = "Event Start"
$ws1.Range("H1").Value = "Event End"
$ws1.Range("I1").Value = "Event Size"
$ws1.Range("J1").Value = "Treated Molecule Count"

# Data updates
$ws1.Range("J2").Value = 88

$ws1.Range("I3").Value = 8600
$ws1.Range("J3").Value = 87

$ws1.Range("J4").Value = 12

$ws1.Range("J5").Value = 117

$ws1.Range("J6").Value = 98

$ws1.Range("J7").Value = 203

$ws1.Range("J8").Value = 43

$ws1.Range("J9").Value = 44

$ws1.Range("H10").Value = 62168640
$ws1.Range("J10").Value = 22

$ws1.Range("J11").Value = 181

$ws1.Range("J12").Value = 75

$ws1.Range("G13").Value = 21982907
$ws1.Range("H13").Value = 20364964
$ws1.Range("J13").Value = 120

# ---------------------------------------------------------------------
# Sheet 2: CNV_and_Aneuploidy_calls
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("CNV_and_Aneuploidy_calls")

$ws2.Range("B1").Value = "Treated Sample Name"
$ws2.Range("F1").Value = "Chromosome"
$ws2.Range("G1").Value = "Event Start"
$ws2.Range("H1").Value = "Event End"
$ws2.Range("I1").Value = "Event Size"
$ws2.Range("J1").Value = "Treated Fractional Copy Number"

# ---------------------------------------------------------------------
# Sheet 3: CNV_metrics
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CNV_metrics")

$ws3.Range("C1").Value = "Treated QC Passed"

$ws3.Range("C2").Value = "yes"
$ws3.Range("E2").Value = "yes"

$ws3.Range("C3").Value = "yes"
$ws3.Range("E3").Value = "yes"

$ws3.Range("C4").Value = "yes"
$ws3.Range("E4").Value = "yes"

$ws3.Range("C5").Value = "yes"
$ws3.Range("E5").Value = "yes"
